$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "YGptT120"
$ws.Range("B2").Value = 23072636
$ws.Range("C2").Value = "yvzdqcq33"
$ws.Range("D2").Value = "Hx4&3n#W"
$ws.Range("F2").Value = "hhvrnRtD"
$ws.Range("G2").Value = "MKqD"
